$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (column G)
$wsOverview.Range("G2").Value = "2016-09-02 08:18:47"
$wsOverview.Range("G5").Value = "2016-09-02 08:18:47"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-02 08:18:42"
$wsZhCn.Range("H5").Value = "2016-09-02 08:18:42"
$wsZhCn.Range("K2").Value = "2016-09-02 08:18:59"
$wsZhCn.Range("K5").Value = "2016-09-02 08:18:59"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-02 08:18:47"
$wsDeDe.Range("H5").Value = "2016-09-02 08:18:47"
$wsDeDe.Range("K2").Value = "2016-09-02 08:19:14"
$wsDeDe.Range("K5").Value = "2016-09-02 08:19:14"
